# feat(sheet): add the function give_names_of_maximum and the corresponding
# command maxnames
#
# The data sheet ("sheet1") holds, for every student row (2..15), a run of
# "group" label cells in columns C..BP followed by a trailing result cell
# (a rank number, or blank). The new command/feature needs 12 extra "group"
# columns, so the run of group-label columns is extended from C:BO (65 cols)
# to C:CA (77 cols), and the trailing result cell that used to live in BP is
# pushed out to CB.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 15
$firstNewCol = 68   # BP (first of the 12 newly used "group" columns)
$lastNewCol  = 79   # CA (last of the 12 newly used "group" columns)
$finalCol    = 80   # CB (new home of the trailing rank/blank cell)
$groupCol    = 67   # BO (pre-existing "group" column, used as the
                     # template value for the newly used columns)

for ($r = 2; $r -le $lastDataRow; $r++) {
    $groupValue = $ws.Cells.Item($r, $groupCol).Value()
    $finalValue = $ws.Cells.Item($r, $firstNewCol).Value()

    $newRange = $ws.Range($ws.Cells.Item($r, $firstNewCol), $ws.Cells.Item($r, $lastNewCol))
    $finalCell = $ws.Cells.Item($r, $finalCol)

    # Touch the formatting (a no-op re-assignment) so every cell in the
    # range materializes in the sheet even when it stays blank (matches
    # rows such as row 10, which has no group label at all).
    $newRange.Font.Name = $newRange.Font.Name
    $finalCell.Font.Name = $finalCell.Font.Name

    if ($groupValue -ne $null -and $groupValue -ne "") {
        $newRange.Value = $groupValue
    }

    if ($finalValue -ne $null -and $finalValue -ne "") {
        $finalCell.Value = $finalValue
    }
}
